$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '87.369.35'
$c.ClearFormats()
$ws.Range("E2").Value = '  -0.56%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.222.87'
$c.ClearFormats()
$ws.Range("E3").Value = '  -3.14%  '

$ws.Range("E4").Value = '  +0.27%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '205.28'
$c.ClearFormats()
$ws.Range("E5").Value = '  -6.63%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '610.31'
$c.ClearFormats()
$ws.Range("E6").Value = '  -6.23%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.377'
$c.ClearFormats()
$ws.Range("E7").Value = '  +5.20%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.675'
$c.ClearFormats()
$ws.Range("E8").Value = '  +11.39%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E9").Value = '  +0.14%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '3.221.95'
$c.ClearFormats()
$ws.Range("E10").Value = '  -3.07%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.542'
$c.ClearFormats()
$ws.Range("E11").Value = '  -7.80%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.180'
$c.ClearFormats()
$ws.Range("E12").Value = '  +6.71%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000247'
$c.ClearFormats()
$ws.Range("E13").Value = '  -9.12%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.825.65'
$c.ClearFormats()
$ws.Range("E14").Value = '  -2.95%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.31'
$c.ClearFormats()
$ws.Range("E15").Value = '  -3.25%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '32.69'
$c.ClearFormats()
$ws.Range("E16").Value = '  -8.27%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '87.337.63'
$c.ClearFormats()
$ws.Range("E17").Value = '  -0.32%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.236.76'
$c.ClearFormats()
$ws.Range("E18").Value = '  -2.53%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.ClearFormats()
$ws.Range("E19").Value = '  -6.05%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.45'
$c.ClearFormats()
$ws.Range("E20").Value = '  -8.95%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '421.90'
$c.ClearFormats()
$ws.Range("E21").Value = '  -7.76%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.57'
$c.ClearFormats()
$ws.Range("E22").Value = '  -13.45%  '

$ws.Range("E23").Value = '  -7.54%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.20'
$c.ClearFormats()
$ws.Range("E24").Value = '  -8.13%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '11.66'
$c.ClearFormats()
$ws.Range("E25").Value = '  -8.52%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.390.21'
$c.ClearFormats()
$ws.Range("E26").Value = '  -2.75%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '74.48'
$c.ClearFormats()
$ws.Range("E27").Value = '  -5.82%  '

$ws.Range("E28").Value = '  +4.63%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -12.90%  '

$ws.Range("E31").Value = '  +0.11%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '545.88'
$c.ClearFormats()
$ws.Range("E32").Value = '  -10.77%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '8.44'
$c.ClearFormats()
$ws.Range("E33").Value = '  -11.02%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.ClearFormats()
$ws.Range("E34").Value = '  -10.81%  '

$ws.Range("E35").Value = '  -20.97%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.69'
$c.ClearFormats()
$ws.Range("E36").Value = '  -7.94%  '

$ws.Range("E37").Value = '  -8.28%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '22.29'
$c.ClearFormats()
$ws.Range("E38").Value = '  -4.56%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '21.83'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.02%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.00'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.71%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.382'
$c.ClearFormats()
$ws.Range("E42").Value = '  -9.51%  '

$ws.Range("E43").Value = '  -0.09%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.90'
$c.ClearFormats()
$ws.Range("E44").Value = '  -12.39%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '146.50'
$c.ClearFormats()
$ws.Range("E45").Value = '  -8.28%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '174.46'
$c.ClearFormats()
$ws.Range("E46").Value = '  -9.07%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '43.72'
$c.ClearFormats()
$ws.Range("E47").Value = '  -5.43%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.129'
$c.ClearFormats()
$ws.Range("E48").Value = '  +12.50%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.28'
$c.ClearFormats()
$ws.Range("E49").Value = '  -10.47%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '4.06'
$c.ClearFormats()
$ws.Range("E50").Value = '  -9.42%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.604'
$c.ClearFormats()
$ws.Range("E51").Value = '  -8.96%  '
